$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove column A entirely; remaining columns (B:F -> A:E) shift one place left.
$ws.Range("A1").EntireColumn.Delete()
